$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: locate the paragraph that contains the first occurrence of
# $searchText at or after document position $afterPos.
# ---------------------------------------------------------------------
function FindParaContaining($searchText, $afterPos) {
    $rng = $d.Range($afterPos, $d.Content.End)
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) { return $null }
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Start -le $rng.Start -and $p.Range.End -ge $rng.End) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# Helper: change a number that appears (as the first occurrence of
# $oldNumber) inside $para into $newNumber, reproducing the same run
# layout Word itself would produce when a user selects the leading
# digit(s) of the old number and retypes them, leaving the final
# digit (and whatever text follows it) untouched:
#   "(-50)"  -> select "5", type "10" -> "(-" / "10" / "0)"
#   "(20)"   -> select "2", type "5"  -> "("  / "5"  / "0)"
# ---------------------------------------------------------------------
function ReplaceNumberAndSplit($para, $oldNumber, $newNumber) {
    $txt = $para.Range.Text
    $idx = $txt.IndexOf($oldNumber)
    $paraStart = $para.Range.Start

    $oldFirstChar = $oldNumber.Substring(0, 1)
    $newPrefix = $newNumber.Substring(0, $newNumber.Length - 1)

    $start = $paraStart + $idx
    $end = $start + $oldFirstChar.Length
    $target = $d.Range($start, $end)
    $target.Text = $newPrefix

    # Force the engine to split the newly-typed text into its own run
    # (mirrors what happens in real Word when a selection is retyped)
    # by toggling a direct-character-formatting property on it and
    # then clearing it again.
    $mid = $d.Range($start, $start + $newPrefix.Length)
    $mid.Bold = $true
    $mid.Bold = $false
}

# --- Locate the three paragraphs inside the "Model #8" block ---------
$modelHeading = FindParaContaining "Model #8" 0
$afterModel8 = $modelHeading.Range.End

$penaltyPara = FindParaContaining "penalty for revisiting the nodes with repetition mask" $afterModel8
$rewardPara  = FindParaContaining "reward for exploring new nodes" $afterModel8

# --- Edit 1: "Model #" / "8" / ": " -> single run "Model #8: " -------
$modelHeading.Range.Find.Execute("Model #8: ", $true, $false, $false, $false, $false, $true, 1, $false, "Model #8: ", 2) | Out-Null

# --- Edit 2: repetition penalty -50 -> -100 ---------------------------
ReplaceNumberAndSplit $penaltyPara "50" "100"

# --- Edit 3: exploration reward 20 -> 50 ------------------------------
ReplaceNumberAndSplit $rewardPara "20" "50"
